$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.520.42"
$ws.Range("E2").Value = "  +4.73%  "

$ws.Range("D3").Value = "1.601.07"
$ws.Range("E3").Value = "  +2.55%  "

$ws.Range("E4").Value = "  -0.45%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.34%  "

$ws.Range("E6").Value = "  +1.70%  "

$ws.Range("E7").Value = "  -0.56%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.16"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.99%  "

$ws.Range("E9").Value = "  +1.59%  "

$ws.Range("E10").Value = "  +1.18%  "

$ws.Range("E11").Value = "  +2.21%  "

$ws.Range("D12").Value = "1.827.61"
$ws.Range("E12").Value = "  +2.15%  "

$ws.Range("D13").Value = "1.596.81"
$ws.Range("E13").Value = "  +1.92%  "

$ws.Range("E14").Value = "  +1.21%  "

$ws.Range("E15").Value = "  +3.61%  "

$ws.Range("D16").Value = "28.523.13"
$ws.Range("E16").Value = "  +4.94%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.62%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.77%  "

$ws.Range("D20").Value = "0.0₃0712"
$ws.Range("E20").Value = "  +1.77%  "

$ws.Range("E21").Value = "  -0.39%  "

$ws.Range("E22").Value = "  +0.39%  "

$ws.Range("E23").Value = "  +2.65%  "

$ws.Range("E24").Value = "  +0.81%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.18%  "

$ws.Range("E26").Value = "  +2.29%  "

$ws.Range("E27").Value = "  +0.28%  "

$ws.Range("E28").Value = "  +1.41%  "

$ws.Range("E29").Value = "  -0.40%  "

$ws.Range("E30").Value = "  +0.95%  "

$ws.Range("E31").Value = "  +1.42%  "

$ws.Range("E32").Value = "  +0.80%  "

$ws.Range("E33").Value = "  +1.18%  "

$ws.Range("D34").Value = "1.424.26"
$ws.Range("E34").Value = "  -0.61%  "

$ws.Range("E35").Value = "  +0.14%  "

$ws.Range("E36").Value = "  -3.63%  "

$ws.Range("E37").Value = "  -0.20%  "

$ws.Range("E38").Value = "  +1.04%  "

$ws.Range("E39").Value = "  +2.94%  "

$ws.Range("E40").Value = "  +8.36%  "

$ws.Range("E41").Value = "  +2.43%  "

$ws.Range("E42").Value = "  -2.77%  "

$ws.Range("E43").Value = "  -0.40%  "

$ws.Range("E44").Value = "  +6.80%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.980"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.23%  "

$ws.Range("D47").Value = "1.739.36"
$ws.Range("E47").Value = "  +2.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.33%  "

$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("E50").Value = "  +5.78%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0526"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.26%  "
